# The commit adds a missing "Primera" quality price record for
# Vega Monumental Concepción / Tomate / Larga vida, dated 2021-09-08
# (Excel serial 44447), which should sit right before the existing
# "Segunda" record for that same date currently on row 326.
#
# Net effect: insert one new row at row 326 (pushing the existing
# rows 326:398 down to 327:399) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 326, shifting
# rows 326:398 down to 327:399.
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new price record.
$ws.Range("A326").Value = 11
$ws.Range("B326").Value = 'Vega Monumental Concepción'
$ws.Range("C326").Value = 'Bíobío'
$ws.Range("D326").Value = 44447
$ws.Range("E326").Value = 8
$ws.Range("F326").Value = 100112020
$ws.Range("G326").Value = 'Tomate'
$ws.Range("H326").Value = 'Larga vida'
$ws.Range("I326").Value = 'Primera'
$ws.Range("J326").Value = 600
$ws.Range("K326").Value = 21000
$ws.Range("L326").Value = 22000
$ws.Range("M326").Value = 21500
$ws.Range("N326").Value = '$/bandeja 18 kilos'
$ws.Range("O326").Value = 'Región de Arica y Parinacota'
$ws.Range("P326").Value = 1194
$ws.Range("Q326").Value = 18
$ws.Range("R326").Value = 'Hortaliza'
